# Fixed all tests, added copy method to Arithmetic_Dict
#
# 1) test_market sheet (sheet1): fix typo "avalilable" -> "available" in R3.
#    This also causes the now-unused "avalilable" shared string to be dropped,
#    shifting every other shared-string index down by one (handled automatically
#    by the engine).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Range("R3").Value = "available"

# 2) test_peasants sheet (sheet2): rows 5 and 16 get the same formatting as the
#    rows above them (yellow fill for labels/ratios, plus the 2-decimal number
#    format used for the "base" values), and new numeric inputs.
#    Row 5 / Row 16 correspond to the "need_ratio" rows of the two example
#    tables.
$ws2.Range("B4").Copy() | Out-Null
$ws2.Range("B5").PasteSpecial(-4122) | Out-Null
$ws2.Range("A4").Copy() | Out-Null
$ws2.Range("A5").PasteSpecial(-4122) | Out-Null
$ws2.Range("A3").Copy() | Out-Null
$ws2.Range("C5").PasteSpecial(-4122) | Out-Null

$ws2.Range("B15").Copy() | Out-Null
$ws2.Range("B16").PasteSpecial(-4122) | Out-Null
$ws2.Range("A15").Copy() | Out-Null
$ws2.Range("A16").PasteSpecial(-4122) | Out-Null
$ws2.Range("A14").Copy() | Out-Null
$ws2.Range("C16").PasteSpecial(-4122) | Out-Null

$ws2.Range("B5").Value = 1.5
$ws2.Range("C5").Value = 1

$ws2.Range("B16").Value = 1.5
$ws2.Range("C16").Value = 1.7

# 3) Restore selections / active sheet state.
#    test_peasants (sheet2) stays the active tab; test_market (sheet1) just
#    gets its selection updated without becoming the active tab.
$ws1.Range("R9").Select() | Out-Null
$ws2.Range("I11").Select() | Out-Null
